$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 167.33333
$ws.Range("I12").Value = 167.33333
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 167.33333
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 2.666670000000011
$ws.Range("N12").ClearContents()
$ws.Range("H31").Value = 1250
$ws.Range("I31").Value = 900
$ws.Range("J31").Value = 3000
$ws.Range("K31").Value = 2700
$ws.Range("L31").Value = 9000
$ws.Range("M31").Value = -2470
$ws.Range("N31").Value = -9460
$ws.Range("H32").Value = 7813300
$ws.Range("I32").Value = 15625632
$ws.Range("J32").Value = 967.75
$ws.Range("K32").Value = 15625632
$ws.Range("L32").Value = 967.75
$ws.Range("M32").Value = -15625306
$ws.Range("N32").Value = -1619.75
$ws.Range("H40").Value = 2591914
$ws.Range("I40").Value = 4466587.5
$ws.Range("J40").Value = 717240.5600000001
$ws.Range("K40").Value = 4466587.5
$ws.Range("L40").Value = 717240.5600000001
$ws.Range("M40").Value = -4466412.5
$ws.Range("N40").Value = -717590.5600000001
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H64").Value = 2542.8572
$ws.Range("I64").Value = 2514.2856
$ws.Range("J64").Value = 2585.7144
$ws.Range("K64").Value = 2514.2856
$ws.Range("L64").Value = 2585.7144
$ws.Range("M64").Value = -2266.2856
$ws.Range("N64").Value = -3081.7144
$ws.Range("H67").Value = 2542.8572
$ws.Range("I67").Value = 2514.2856
$ws.Range("J67").Value = 2585.7144
$ws.Range("K67").Value = 2514.2856
$ws.Range("L67").Value = 2585.7144
$ws.Range("M67").Value = -1656.2856
$ws.Range("N67").Value = -4301.7144
$ws.Range("H76").Value = 55574532
$ws.Range("I76").Value = 52000
$ws.Range("K76").Value = 52000
$ws.Range("M76").Value = -51685
$ws.Range("H79").Value = 55574532
$ws.Range("I79").Value = 52000
$ws.Range("K79").Value = 52000
$ws.Range("M79").Value = -50908
$ws.Range("H92").Value = 607.9091
$ws.Range("J92").Value = 470
$ws.Range("L92").Value = 470
$ws.Range("N92").Value = -2966
$ws.Range("H96").Value = 1358.6666
$ws.Range("J96").Value = 1632.5714
$ws.Range("L96").Value = 4897.7142
$ws.Range("N96").Value = -7643.7142
$ws.Range("H97").Value = 77309460
$ws.Range("J97").Value = 77309460
$ws.Range("L97").Value = 231928380
$ws.Range("N97").Value = -231929372
$ws.Range("H101").Value = 6493956
$ws.Range("I101").Value = 249.33333
$ws.Range("J101").Value = 11364236
$ws.Range("K101").Value = 747.99999
$ws.Range("L101").Value = 34092708
$ws.Range("M101").Value = 874.00001
$ws.Range("N101").Value = -34095952
$ws.Range("H116").Value = 2667.6667
$ws.Range("I116").Value = 2740
$ws.Range("J116").Value = 2623.1538
$ws.Range("K116").Value = 2740
$ws.Range("L116").Value = 2623.1538
$ws.Range("M116").Value = 702
$ws.Range("N116").Value = -9507.1538
$ws.Range("H121").Value = 700
$ws.Range("I121").Value = 400
$ws.Range("J121").Value = 1000
$ws.Range("K121").Value = 1200
$ws.Range("L121").Value = 3000
$ws.Range("M121").Value = 547
$ws.Range("N121").Value = -6494
$ws.Range("H131").Value = 1663.069
$ws.Range("I131").Value = 465.16666
$ws.Range("J131").Value = 7413
$ws.Range("K131").Value = 1395.49998
$ws.Range("L131").Value = 22239
$ws.Range("M131").Value = 3644.50002
$ws.Range("N131").Value = -32319

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 204.44444
$ws.Range("I5").Value = 98.5
$ws.Range("K5").Value = 98.5
$ws.Range("M5").Value = 13.5
$ws.Range("H45").Value = 144421
$ws.Range("I45").Value = 334482.66
$ws.Range("K45").Value = 334482.66
$ws.Range("M45").Value = -334105.66
$ws.Range("H54").Value = 12926.667
$ws.Range("J54").Value = 12926.667
$ws.Range("L54").Value = 12926.667
$ws.Range("N54").Value = -14464.667
$ws.Range("H102").Value = 1252.5
$ws.Range("I102").Value = 999.6667
$ws.Range("K102").Value = 999.6667
$ws.Range("M102").Value = 622.3333
$ws.Range("H122").Value = 1613.1578
$ws.Range("I122").Value = 1434.6154
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 4303.8462
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -1853.8462
$ws.Range("N122").Value = -10900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 204.44444
$ws.Range("I4").Value = 98.5
$ws.Range("K4").Value = 98.5
$ws.Range("M4").Value = 16.5
$ws.Range("H45").Value = 12065
$ws.Range("J45").Value = 12065
$ws.Range("L45").Value = 12065
$ws.Range("N45").Value = -13681
$ws.Range("H86").Value = 665788
$ws.Range("I86").Value = 1179
$ws.Range("J86").Value = 3878064.8
$ws.Range("K86").Value = 1179
$ws.Range("L86").Value = 3878064.8
$ws.Range("M86").Value = -56
$ws.Range("N86").Value = -3880310.8
$ws.Range("H89").Value = 665788
$ws.Range("I89").Value = 1179
$ws.Range("J89").Value = 3878064.8
$ws.Range("K89").Value = 5895
$ws.Range("L89").Value = 19390324
$ws.Range("M89").Value = -279
$ws.Range("N89").Value = -19401556
$ws.Range("H105").Value = 62501636
$ws.Range("I105").Value = 1680
$ws.Range("J105").Value = 250001500
$ws.Range("K105").Value = 1680
$ws.Range("L105").Value = 250001500
$ws.Range("M105").Value = 67
$ws.Range("N105").Value = -250004994

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 11483.333
$ws.Range("I7").Value = 349.33334
$ws.Range("K7").Value = 349.33334
$ws.Range("M7").Value = -236.33334
$ws.Range("H105").Value = 9806.392
$ws.Range("I105").Value = 11145.632
$ws.Range("K105").Value = 11145.632
$ws.Range("M105").Value = -9398.632
$ws.Range("H122").Value = 25000504
$ws.Range("I122").Value = 50000384
$ws.Range("J122").Value = 624
$ws.Range("K122").Value = 150001152
$ws.Range("L122").Value = 1872
$ws.Range("M122").Value = -149998702
$ws.Range("N122").Value = -6772

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 283.33334
$ws.Range("I23").Value = 200
$ws.Range("J23").Value = 300
$ws.Range("K23").Value = 600
$ws.Range("L23").Value = 900
$ws.Range("M23").Value = -365
$ws.Range("N23").Value = -1370
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H109").Value = 3890.3125
$ws.Range("I109").Value = 1417.6666
$ws.Range("J109").Value = 4460.923
$ws.Range("K109").Value = 4252.9998
$ws.Range("L109").Value = 13382.769
$ws.Range("M109").Value = -3212.9998
$ws.Range("N109").Value = -15462.769
$ws.Range("H132").Value = 22733058
$ws.Range("I132").Value = 775.3333
$ws.Range("J132").Value = 38470790
$ws.Range("K132").Value = 6977.9997
$ws.Range("L132").Value = 346237110
$ws.Range("M132").Value = -4447.9997
$ws.Range("N132").Value = -346242170

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 19642.572
$ws.Range("I122").Value = 24833.428
$ws.Range("J122").Value = 4070
$ws.Range("K122").Value = 74500.284
$ws.Range("L122").Value = 12210
$ws.Range("M122").Value = -72050.284
$ws.Range("N122").Value = -17110
